# Adding Button definition class
#
# Insert a new sub-bullet ("Create UI before implementing ECS") right
# after the "Create a System base class or concept." bullet and before
# the "Implement MovementSystem:" bullet, under Stage 2: Core Systems.

$d = $word.ActiveDocument

# Locate the "Implement MovementSystem:" paragraph by scanning - the new
# bullet is inserted immediately before it (i.e. right after the
# "...System base class or concept." bullet).
$targetIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs.Item($i)
    if ($p.Range.Text -like "Implement MovementSystem*") {
        $targetIndex = $i
        break
    }
}

$target = $d.Paragraphs.Item($targetIndex)
$insertionPoint = $target.Range
$insertionPoint.Collapse(1)
$insertionPoint.InsertParagraphBefore()

# Re-fetch by index: the newly created (empty) paragraph now occupies the
# slot the target used to be at, inheriting $target's list/numbering/font
# formatting (no strikethrough), while $target itself shifted down by one.
$newPara = $d.Paragraphs.Item($targetIndex)
$newPara.Range.Text = "Create UI before implementing ECS"

Write-Output "Inserted new bullet: $($newPara.Range.Text)"
